$d = $word.ActiveDocument

# 1) Renumber "1.2 Guia de cabo" -> "1.4 Guia de cabo"
$d.Content.Find.Execute("1.2 Guia de cabo", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "1.4 Guia de cabo", 2) | Out-Null

# 2) Retitle the first heading: "1.1 ELETRODUTO PVC RÍGIDO DE ½” ANTICHAMA"
#    becomes "1.1 CAIXA DE PASSAGEM DE SOBREPOR 120x120x7,5CM"
$d.Content.Find.Execute("1.1 ELETRODUTO PVC RÍGIDO DE ½” ANTICHAMA", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "1.1 CAIXA DE PASSAGEM DE SOBREPOR 120x120x7,5CM", 2) | Out-Null

# 3) Locate that (now renamed) heading paragraph so we can splice the new
#    sections in right after it, before the pre-existing "Eletroduto
#    fabricado..." bullet list.
$headingIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "1.1 CAIXA DE PASSAGEM DE SOBREPOR 120x120x7,5CM") {
        $headingIndex = $i
        break
    }
}

if ($headingIndex -eq -1) {
    throw "Could not locate the '1.1 CAIXA DE PASSAGEM DE SOBREPOR 120x120x7,5CM' heading paragraph."
}

$heading = $d.Paragraphs.Item($headingIndex)
$heading.Range.InsertParagraphAfter() | Out-Null
$newPara = $d.Paragraphs.Item($headingIndex + 1)

$xml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:pPr><w:pStyle w:val="ListBullet"/></w:pPr><w:r><w:t>Caixa de sobrepor na cor cinza ou bege com tampa fixada por pressão;</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="ListBullet"/></w:pPr><w:r><w:t>Possuir grau de proteção IP44;</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="ListBullet"/></w:pPr><w:r><w:t>Deverá possuir marcação para entrada de eletrodutos e canaletas;</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="ListBullet"/></w:pPr><w:r><w:t>Deverá possuir dimensões de no mínimo 120x120x75mm;</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="ListBullet"/></w:pPr><w:r><w:t>Deverá ser fornecido com arruelas para a conexão dos eletrodutos conforme necessidade.</w:t></w:r></w:p>
<w:p><w:r/></w:p>
<w:p><w:r><w:rPr><w:b/><w:color w:val="4F81BD"/></w:rPr><w:t>1.2 ELETRODUTO GALVANIZADO A FOGO DE 1"</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="ListBullet"/></w:pPr><w:r><w:t>Deverá ser fornecido com diâmetro de 1”;</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="ListBullet"/></w:pPr><w:r><w:t>Deverá ser fabricado em aço com acabamento galvanizado a fogo (imersão a quente) e indicado para instalação em ambientes externos;</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="ListBullet"/></w:pPr><w:r><w:t>Deverá ser fornecido com todos os acessórios necessários para sua instalação tais como, curva, luva, abraçadeira, parafusos para fixação, tirante rosqueado, cantoneiras etc.</w:t></w:r></w:p>
<w:p><w:r/></w:p>
<w:p><w:r><w:rPr><w:b/><w:color w:val="4F81BD"/></w:rPr><w:t>1.3 ELETRODUTO PVC RÍGIDO DE ½” ANTICHAMA</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$newPara.Range.InsertXML($xml) | Out-Null

Write-Output "done"
